# "Partial Script for add revision"
# Updates the ExecutionFlag column on the ModuleController sheet:
#   - Login   (row 2): Yes -> No
#   - CommonTC(row 6): No  -> Yes
#   - SMOKE   (row 7): NO  -> Yes   (also drops the stray "NO" shared string)
# and moves the saved cell selection from I17 to H16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModuleController")

$ws.Range("A2").Value = "No"
$ws.Range("A6").Value = "Yes"
$ws.Range("A7").Value = "Yes"

# Match the saved selection/active cell recorded in the workbook view.
$ws.Range("H16").Select()

# Best-effort: nudge the window's vertical screen position (cosmetic).
try {
    $excel.ActiveWindow.Top = 4140
} catch {
}
